$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the dates between rows 2/3 (were 2022-01-13) and rows 6/7 (were 2021-12-29)
$ws.Range("D2").Value = 44559
$ws.Range("D3").Value = 44559
$ws.Range("D6").Value = 44574
$ws.Range("D7").Value = 44574
